$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 values (C7/D7/E7 look numeric, e.g. "23456") must be stored as
# literal text, matching the rest of the sheet (all cells are text).
# Force text entry by setting the number format to Text ("@") before
# writing, then clear the formatting again so the new row keeps the
# sheet's plain default styling.
$ws.Range("A7:E7").NumberFormat = "@"

$ws.Range("A7").Value = "chitti"
$ws.Range("B7").Value = "chitti@gmail.com"
$ws.Range("C7").Value = "23456"
$ws.Range("D7").Value = "23"
$ws.Range("E7").Value = "12345"

$ws.Range("A7:E7").ClearFormats()
